$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2024-05-30 Thursday" "2024-05-31 Friday"

Replace-Text "65÷7=9, 2" "44÷7=6, 2"
Replace-Text "73÷2=36, 1" "70÷7=10, 0"
Replace-Text "98÷8=12, 2" "45÷4=11, 1"
Replace-Text "71÷6=11, 5" "47÷3=15, 2"
Replace-Text "88÷7=12, 4" "66÷6=11, 0"

Replace-Text "37÷8=4, 5" "89÷4=22, 1"
Replace-Text "33÷9=3, 6" "42÷5=8, 2"
Replace-Text "85÷9=9, 4" "35÷3=11, 2"
Replace-Text "58÷5=11, 3" "98÷9=10, 8"
Replace-Text "47÷7=6, 5" "65÷4=16, 1"

Replace-Text "90÷9=10, 0" "23÷4=5, 3"
Replace-Text "29÷8=3, 5" "48÷4=12, 0"
Replace-Text "15÷8=1, 7" "53÷9=5, 8"
Replace-Text "40÷6=6, 4" "76÷9=8, 4"
Replace-Text "78÷5=15, 3" "73÷3=24, 1"

Replace-Text "69÷8=8, 5" "36÷8=4, 4"
Replace-Text "53÷7=7, 4" "71÷5=14, 1"
Replace-Text "32÷5=6, 2" "92÷6=15, 2"
Replace-Text "47÷2=23, 1" "14÷3=4, 2"
Replace-Text "18÷6=3, 0" "50÷5=10, 0"

Replace-Text "41÷5=8, 1" "57÷9=6, 3"
Replace-Text "61÷9=6, 7" "35÷7=5, 0"
Replace-Text "85÷5=17, 0" "11÷3=3, 2"
Replace-Text "43÷4=10, 3" "35÷9=3, 8"
Replace-Text "17÷3=5, 2" "77÷6=12, 5"
